$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '70.429.65'
$ws.Range('E2').Value = '  -2.23%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.624.39'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.18'
$ws.Range('E5').Value = '  -2.53%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.32'
$ws.Range('E6').Value = '  -4.07%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.615.95'
$ws.Range('E7').Value = '  +0.15%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('E9').Value = '  +0.09%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.196'
$ws.Range('E10').Value = '  -5.23%  '

$ws.Range('E11').Value = '  +17.10%  '

$ws.Range('E12').Value = '  -0.85%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '48.32'
$ws.Range('E13').Value = '  -4.34%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000283'
$ws.Range('E14').Value = '  -2.33%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.210.36'
$ws.Range('E15').Value = '  +0.26%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '674.46'
$ws.Range('E16').Value = '  -4.54%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '8.90'
$ws.Range('E17').Value = '  -0.73%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.623.96'
$ws.Range('E18').Value = '  +1.61%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.475.65'
$ws.Range('E19').Value = '  -2.09%  '

$ws.Range('E20').Value = '  -0.46%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.76'
$ws.Range('E21').Value = '  -4.27%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.45'
$ws.Range('E22').Value = '  -2.79%  '

$ws.Range('E23').Value = '  +0.29%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '17.03'
$ws.Range('E24').Value = '  -4.69%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '99.55'
$ws.Range('E25').Value = '  -5.85%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.91'
$ws.Range('E26').Value = '  -3.08%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.79'
$ws.Range('E27').Value = '  -2.68%  '

$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('E29').Value = '  -1.90%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '34.57'
$ws.Range('E30').Value = '  -3.34%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.02'
$ws.Range('E31').Value = '  -1.11%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.29'
$ws.Range('E32').Value = '  -4.60%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.56'
$ws.Range('E33').Value = '  +1.52%  '

$ws.Range('E34').Value = '  -7.13%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.00'
$ws.Range('E35').Value = '  -4.80%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '576.44'
$ws.Range('E36').Value = '  -2.66%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '11.03'
$ws.Range('E37').Value = '  -3.30%  '

$ws.Range('E38').Value = '  -0.98%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '58.16'
$ws.Range('E39').Value = '  -2.84%  '

$ws.Range('E40').Value = '  +0.18%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.559.47'
$ws.Range('E41').Value = '  -2.65%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.141'
$ws.Range('E42').Value = '  -2.97%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0450'
$ws.Range('E43').Value = '  -0.80%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.344'
$ws.Range('E44').Value = '  -1.16%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '34.40'
$ws.Range('E45').Value = '  -4.66%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₃0728'
$ws.Range('E46').Value = '  -6.17%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.68'
$ws.Range('E47').Value = '  -4.42%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.84'
$ws.Range('E48').Value = '  +1.81%  '

$ws.Range('E49').Value = '  -0.19%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '136.52'
$ws.Range('E50').Value = '  +2.31%  '

$ws.Range('E51').Value = '  -1.86%  '
